# Generate Report for Handback
# Populates the "Latest Target File" (I), "Latest Handback File" (J) and
# "Latest Handback DateTime" (K) columns for the 861fd803-... row, now that a
# handback was produced (but flagged as stale), on both the zh-cn and de-de
# status sheets, and records the "not the latest" error detail in column P.

$wb = $excel.ActiveWorkbook

$handbackMdDisplay = "861fd803-29f5-4fd6-b602-079be2af0ede.md"
$handbackMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d70bee1ba72206d25656ac4d26c790dd3acef31/e2e/861fd803-29f5-4fd6-b602-079be2af0ede.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03f2765267d0347a45c5208d8cdefa8ab59098d6/e2e/861fd803-29f5-4fd6-b602-079be2af0ede.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d70bee1ba72206d25656ac4d26c790dd3acef31/e2e/861fd803-29f5-4fd6-b602-079be2af0ede.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7 (861fd803-29f5-4fd6-b602-079be2af0ede)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("K7").Value = "2016-08-24 06:56:09"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Range("J7").Value = $wsZh.Range("G7").Value2

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackMdUrl, "", "", $handbackMdDisplay)

# ---------------------------------------------------------------------
# de-de sheet, row 7 (861fd803-29f5-4fd6-b602-079be2af0ede)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("K7").Value = "2016-08-24 06:56:16"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Range("J7").Value = $wsDe.Range("G7").Value2

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackMdUrl, "", "", $handbackMdDisplay)
